$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header formatting from an existing header cell (H1) to the new
# header cells I1 and J1 so they get the same bold/centered/bordered style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Set the new header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Set the new data values in row 2
$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 8
